$wb = $excel.ActiveWorkbook

# "OFF" sheet - row 2 (H) target depth data update for Week 13
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 425
$wsOff.Range("C2").Value = 273
$wsOff.Range("D2").Value = 124
$wsOff.Range("E2").Value = 52
$wsOff.Range("F2").Value = 20
$wsOff.Range("G2").Value = 11

# "DEF" sheet - row 2 (H) target depth data update for Week 13
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 407
$wsDef.Range("C2").Value = 275
$wsDef.Range("D2").Value = 121
$wsDef.Range("E2").Value = 54
$wsDef.Range("F2").Value = 12
$wsDef.Range("G2").Value = 10
